$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-5.3(se=1)"
$ws.Range("C2").Value = "0.0011(95% CI, 0.00069-0.0017)"
$ws.Range("B3").Value = "-5.5(se=1.5)"
$ws.Range("B4").Value = "-16.5(se=10.8)"
$ws.Range("B5").Value = "-19.8(se=14.2)"
$ws.Range("G5").Value = "0.015(95% CI, 0.0011-0.067)"
